$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append three new rows of data (141, 142, 143) to the table
$ws.Range("A142").Value = 141
$ws.Range("B142").Value = "Архипов Станислав "
$ws.Range("C142").Value = "Наука"

$ws.Range("A143").Value = 142
$ws.Range("B143").Value = "Струна Марина "
$ws.Range("C143").Value = "Наука"

$ws.Range("A144").Value = 143
$ws.Range("B144").Value = "Седов Станислав"
$ws.Range("C144").Value = "Наука"

# Match the scrolled/selected view left behind after data entry
$ws.Range("B143").Select()
